# Update crypto price/volume figures per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.451.51'
$ws.Range("E2").Value = '  -1.84%  '
$ws.Range("D3").Value = '2.334.91'
$ws.Range("E3").Value = '  -4.47%  '
$ws.Range("D5").Formula = "'541.46"
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").Formula = "'135.46"
$ws.Range("E6").Value = '  -7.07%  '
$ws.Range("D7").Formula = "'0.998"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  -10.91%  '
$ws.Range("E9").Value = '  -4.46%  '
$ws.Range("E10").Value = '  -2.44%  '
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("E12").Value = '  -3.48%  '
$ws.Range("E13").Value = '  -2.94%  '
$ws.Range("D14").Formula = "'24.32"
$ws.Range("E14").Value = '  -5.96%  '
$ws.Range("E15").Value = '  -4.64%  '
$ws.Range("D16").Value = '60.166.32'
$ws.Range("E16").Value = '  -1.94%  '
$ws.Range("E17").Value = '  -5.37%  '
$ws.Range("E18").Value = '  -4.45%  '
$ws.Range("E19").Value = '  -2.09%  '
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("E22").Value = '  -5.72%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Formula = "'62.62"
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("D25").Formula = "'1.66"
$ws.Range("E25").Value = '  -10.24%  '
$ws.Range("E26").Value = '  +6.11%  '
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("E28").Value = '  -4.77%  '
$ws.Range("E29").Value = '  -4.19%  '
$ws.Range("E30").Value = '  -8.34%  '
$ws.Range("D31").Formula = "'488.81"
$ws.Range("E31").Value = '  -6.83%  '
$ws.Range("E32").Value = '  -12.62%  '
$ws.Range("E33").Value = '  -2.32%  '
$ws.Range("E34").Value = '  -4.96%  '
$ws.Range("E35").Value = '  -6.44%  '
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("E37").Value = '  -4.41%  '
$ws.Range("D38").Formula = "'18.39"
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("E40").Value = '  -8.85%  '
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("D42").Formula = "'141.27"
$ws.Range("E42").Value = '  +1.92%  '
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("D44").Formula = "'40.39"
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").Formula = "'140.18"
$ws.Range("E45").Value = '  -1.63%  '
$ws.Range("E46").Value = '  -2.25%  '
$ws.Range("E47").Value = '  -11.00%  '
$ws.Range("E48").Value = '  -3.32%  '
$ws.Range("E49").Value = '  -10.72%  '
$ws.Range("E50").Value = '  -3.90%  '
$ws.Range("E51").Value = '  -4.00%  '
